# Update "想去人数" (F column) counts across the four sheets to match the
# refreshed data snapshot (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1319
$ws1.Range("F4").Value = 207
$ws1.Range("F5").Value = 1032
$ws1.Range("F7").Value = 179
$ws1.Range("F8").Value = 771
$ws1.Range("F9").Value = 1468
$ws1.Range("F10").Value = 1053
$ws1.Range("F11").Value = 789
$ws1.Range("F12").Value = 40467
$ws1.Range("F13").Value = 801
$ws1.Range("F14").Value = 93
$ws1.Range("F15").Value = 606
$ws1.Range("F16").Value = 114
$ws1.Range("F17").Value = 699
$ws1.Range("F18").Value = 1319
$ws1.Range("F19").Value = 212
$ws1.Range("F23").Value = 5361
$ws1.Range("F24").Value = 293
$ws1.Range("F26").Value = 2485
$ws1.Range("F27").Value = 5950
$ws1.Range("F29").Value = 1035
$ws1.Range("F30").Value = 613
$ws1.Range("F33").Value = 1064
$ws1.Range("F35").Value = 80
$ws1.Range("F37").Value = 731
$ws1.Range("F38").Value = 18
$ws1.Range("F41").Value = 1089
$ws1.Range("F44").Value = 58
$ws1.Range("F46").Value = 105
$ws1.Range("F47").Value = 590

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 2095
$ws2.Range("F8").Value = 106
$ws2.Range("F10").Value = 143
$ws2.Range("F12").Value = 686
$ws2.Range("F25").Value = 532
$ws2.Range("F37").Value = 196
$ws2.Range("F39").Value = 505
$ws2.Range("F41").Value = 39
$ws2.Range("F44").Value = 80
$ws2.Range("F45").Value = 102

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 0
$ws3.Range("F4").Value = 677
$ws3.Range("F5").Value = 785
$ws3.Range("F6").Value = 441

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1319
$ws4.Range("F7").Value = 441
$ws4.Range("F10").Value = 207
$ws4.Range("F12").Value = 179
$ws4.Range("F14").Value = 1468
$ws4.Range("F16").Value = 1053
$ws4.Range("F17").Value = 789
$ws4.Range("F18").Value = 93
$ws4.Range("F19").Value = 1319
$ws4.Range("F20").Value = 212
$ws4.Range("F24").Value = 293
$ws4.Range("F25").Value = 2485
$ws4.Range("F26").Value = 5950
$ws4.Range("F28").Value = 1035
$ws4.Range("F31").Value = 613
$ws4.Range("F33").Value = 1064
$ws4.Range("F34").Value = 80
$ws4.Range("F36").Value = 731
$ws4.Range("F39").Value = 1089
$ws4.Range("F43").Value = 39
$ws4.Range("F44").Value = 105
$ws4.Range("F46").Value = 80
$ws4.Range("F47").Value = 102
